$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the numeric cells and the label cell.
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Format B1 with bold font, thin box border, centered/top aligned.
$cell = $ws.Range("B1")
$cell.Font.Bold = $true
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4160
$cell.Borders.LineStyle = 1
$cell.Borders.Weight = 2

# Copy that same formatting onto A2 so both cells share the single
# resulting cell style instead of generating extra intermediate styles.
$cell.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
